$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap column B <-> column F (header + row2/row3 values) ---
$bHeader = $ws.Range("B1").Value()
$fHeader = $ws.Range("F1").Value()
$ws.Range("B1").Value = $fHeader
$ws.Range("F1").Value = $bHeader

$b2 = $ws.Range("B2").Value()
$f2 = $ws.Range("F2").Value()
$ws.Range("B2").Value = $f2
$ws.Range("F2").Value = $b2

$b3 = $ws.Range("B3").Value()
$f3 = $ws.Range("F3").Value()
$ws.Range("B3").Value = $f3
$ws.Range("F3").Value = $b3

# --- Swap column C <-> column G (header + row2/row3 values) ---
$cHeader = $ws.Range("C1").Value()
$gHeader = $ws.Range("G1").Value()
$ws.Range("C1").Value = $gHeader
$ws.Range("G1").Value = $cHeader

$c2 = $ws.Range("C2").Value()
$g2 = $ws.Range("G2").Value()
$ws.Range("C2").Value = $g2
$ws.Range("G2").Value = $c2

$c3 = $ws.Range("C3").Value()
$g3 = $ws.Range("G3").Value()
$ws.Range("C3").Value = $g3
$ws.Range("G3").Value = $c3

# --- Tiny re-computed values for columns D and E (model recalculation) ---
$ws.Range("D2").Value = 31.61408023313981
$ws.Range("E2").Value = 2.648915369657066
$ws.Range("D3").Value = 0.003703943600203319
$ws.Range("E3").Value = 0
